$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new question row (row 9) - bug fix: first matching question at
# requested level is now returned, so this new question row was added back.
$ws.Range("A9").Value = "Pro kola je typycký/á/é:"
$ws.Range("B9").Value = "bradka"
$ws.Range("C9").Value = "hříva"
$ws.Range("D9").Value = "ploutev"
$ws.Range("E9").Value = "tesáky"
$ws.Range("F9").Value = 1

# Update the active selection to reflect where the cursor ended up after editing
[void]$ws.Range("L10").Select()
